$wb = $excel.ActiveWorkbook

$wsRecette = $wb.Worksheets.Item("Recette Client")
$wsLivrables = $wb.Worksheets.Item("Livrables")

# --- Update "Livrables" sheet (descriptif fonctionnel / avancement) ---

# Row 16 -> liv_15 "Capture d'image de l'aspect du site sur tablette, smartphone et Portable"
# Avancement 0 -> 0.9, add comment about the Paypal button / date field fix on step E1
$wsLivrables.Range("C16").Value = 0.9
$wsLivrables.Range("D16").Value = "Action de correction de rendu Bouton Paypal et champ date de l'étape E1"

# Row 17 -> liv_16 "Descriptif fonctionnel du produit"
# Avancement 0 -> 0.7
$wsLivrables.Range("C17").Value = 0.7

# --- Update selections / active cells to match latest user navigation ---
# Select the relevant cell on "Livrables" first (keeps its own stored selection state)
$wsLivrables.Range("D17").Select()

# Finish with "Recette Client" active/selected, as it was before the edit
$wsRecette.Activate()
$wsRecette.Range("D3").Select()
